$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# --- F2: update rpc-reply message-id uuid ---
$f2 = $ws.Range("F2").Value()
$f2 = $f2.Replace(
    'urn:uuid:7e5feaf7-b399-4405-8bff-07a8f0eb629b',
    'urn:uuid:25adee92-df68-4114-b51a-2400ad1e30ab'
)
$ws.Range("F2").Value = $f2

# --- G2: update protocol identifier/name and add peer-as ---
$g2 = $ws.Range("G2").Value()

$g2 = $g2.Replace(
    "<identifier>BGP</identifier>`n              <name>BGP_65100</name>",
    "<identifier xmlns:oc-pol-types=`"http://openconfig.net/yang/policy-types`">oc-pol-types:BGP</identifier>`n              <name>default</name>"
)

$g2 = $g2.Replace(
    "<identifier xmlns:oc-pol-types=`"http://openconfig.net/yang/policy-types`">oc-pol-types:BGP</identifier>`n                <name>BGP_65100</name>",
    "<identifier xmlns:oc-pol-types=`"http://openconfig.net/yang/policy-types`">oc-pol-types:BGP</identifier>`n                <name>default</name>"
)

$g2 = $g2.Replace(
    "<neighbor-address>192.168.1.2</neighbor-address>`n                    </config>",
    "<neighbor-address>192.168.1.2</neighbor-address>`n                      <peer-as>65100</peer-as>`n                    </config>"
)

$ws.Range("G2").Value = $g2
